# Duplicate the existing "email_verification_username_pos" sheet, placing the
# copy BEFORE it, rename the copy to "email_veri_username_neg", update the
# summary numbers on both sheets, set the selection to F7 on both, and leave
# the original sheet ("email_verification_username_pos") as the active tab.

$wb = $excel.ActiveWorkbook

$wsOrig = $wb.Worksheets.Item(1)
$wsOrig.Copy($wsOrig)

$wsNew = $wb.Worksheets.Item(1)
$wsOrig = $wb.Worksheets.Item(2)

$wsNew.Name = "email_veri_username_neg"

$wsNew.Range("F4").Value = 33
$wsNew.Range("F5").Value = 33
$wsNew.Range("F7").Value = 0
$wsNew.Range("F7").Select()

$wsOrig.Range("F4").Value = 26
$wsOrig.Range("F5").Value = 26
$wsOrig.Range("F7").Value = 0
$wsOrig.Range("F7").Select()

$wsOrig.Activate()
